{"js": "// Apply the dated worksheet update: new title date and 25 refreshed\n// two-digit-by-two-digit multiplication problems (same cell positions).\nconst replacements = [\n  [\"2025-12-20 Saturday\", \"2025-12-21 Sunday\"],\n\n  [\"19\u00d735=\", \"90\u00d748=\"],\n  [\"87\u00d725=\", \"76\u00d731=\"],\n  [\"75\u00d783=\", \"26\u00d726=\"],\n  [\"39\u00d764=\", \"23\u00d799=\"],\n  [\"12\u00d786=\", \"58\u00d791=\"],\n\n  [\"36\u00d785=\", \"28\u00d775=\"],\n  [\"38\u00d716=\", \"74\u00d724=\"],\n  [\"59\u00d718=\", \"88\u00d741=\"],\n  [\"43\u00d714=\", \"92\u00d730=\"],\n  [\"55\u00d718=\", \"41\u00d757=\"],\n\n  [\"52\u00d772=\", \"48\u00d734=\"],\n  [\"82\u00d740=\", \"50\u00d765=\"],\n  [\"86\u00d731=\", \"64\u00d770=\"],\n  [\"48\u00d734=\", \"22\u00d753=\"],\n  [\"11\u00d734=\", \"12\u00d787=\"],\n\n  [\"69\u00d757=\", \"47\u00d758=\"],\n  [\"85\u00d757=\", \"63\u00d798=\"],\n  [\"99\u00d713=\", \"72\u00d719=\"],\n  [\"28\u00d771=\", \"11\u00d713=\"],\n  [\"72\u00d771=\", \"69\u00d763=\"],\n\n  [\"75\u00d778=\", \"34\u00d793=\"],\n  [\"70\u00d769=\", \"64\u00d727=\"],\n  [\"54\u00d765=\", \"99\u00d781=\"],\n  [\"66\u00d765=\", \"96\u00d773=\"],\n  [\"87\u00d792=\", \"14\u00d723=\"],\n];\n\nconst body = context.document.body;\n\n// Two-phase replace: some new values coincidentally equal other entries'\n// old values (e.g. \"48\u00d734=\" is both an old and a new value), so a naive\n// single-pass search/replace could re-match a value we just inserted.\n// First swap every old value for a unique placeholder, then swap every\n// placeholder for its real new value.\nconst placeholders = replacements.map((_, i) => `\\u0001PLACEHOLDER_${i}\\u0001`);\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [oldText] = replacements[i];\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(placeholders[i], Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, newText] = replacements[i];\n  const results = body.search(placeholders[i], { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the dated worksheet update: new title date and 25 refreshed\n# two-digit-by-two-digit multiplication problems (same cell positions).\n\n$wdReplaceAll   = 2\n$wdFindContinue = 1\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2025-12-20 Saturday\", \"2025-12-21 Sunday\"),\n\n  @(\"19\u00d735=\", \"90\u00d748=\"),\n  @(\"87\u00d725=\", \"76\u00d731=\"),\n  @(\"75\u00d783=\", \"26\u00d726=\"),\n  @(\"39\u00d764=\", \"23\u00d799=\"),\n  @(\"12\u00d786=\", \"58\u00d791=\"),\n\n  @(\"36\u00d785=\", \"28\u00d775=\"),\n  @(\"38\u00d716=\", \"74\u00d724=\"),\n  @(\"59\u00d718=\", \"88\u00d741=\"),\n  @(\"43\u00d714=\", \"92\u00d730=\"),\n  @(\"55\u00d718=\", \"41\u00d757=\"),\n\n  @(\"52\u00d772=\", \"48\u00d734=\"),\n  @(\"82\u00d740=\", \"50\u00d765=\"),\n  @(\"86\u00d731=\", \"64\u00d770=\"),\n  @(\"48\u00d734=\", \"22\u00d753=\"),\n  @(\"11\u00d734=\", \"12\u00d787=\"),\n\n  @(\"69\u00d757=\", \"47\u00d758=\"),\n  @(\"85\u00d757=\", \"63\u00d798=\"),\n  @(\"99\u00d713=\", \"72\u00d719=\"),\n  @(\"28\u00d771=\", \"11\u00d713=\"),\n  @(\"72\u00d771=\", \"69\u00d763=\"),\n\n  @(\"75\u00d778=\", \"34\u00d793=\"),\n  @(\"70\u00d769=\", \"64\u00d727=\"),\n  @(\"54\u00d765=\", \"99\u00d781=\"),\n  @(\"66\u00d765=\", \"96\u00d773=\"),\n  @(\"87\u00d792=\", \"14\u00d723=\")\n)\n\n# Two-phase replace: some new values coincidentally equal other entries'\n# old values (e.g. \"48\u00d734=\" is both an old and a new value), so a naive\n# single-pass find/replace could re-match text we just inserted. First\n# swap every old value for a unique placeholder, then swap every\n# placeholder for its real new value.\nfor ($i = 0; $i -lt $replacements.Count; $i++) {\n  $oldText = $replacements[$i][0]\n  $placeholder = [string]::Format(\"{0}PLACEHOLDER_{1}{0}\", [char]1, $i)\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $placeholder, $wdReplaceAll) | Out-Null\n}\n\nfor ($i = 0; $i -lt $replacements.Count; $i++) {\n  $newText = $replacements[$i][1]\n  $placeholder = [string]::Format(\"{0}PLACEHOLDER_{1}{0}\", [char]1, $i)\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Execute($placeholder, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll) | Out-Null\n}\n"}
